$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 395, shifting existing rows 395:456 down to 396:457.
# Excel copies formatting (incl. number format) from the row above automatically.
$ws.Rows.Item(395).Insert()

# Populate the new row 395 with its data.
$ws.Range("A395").Value = 7
$ws.Range("B395").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C395").Value = "Ñuble"
$ws.Range("D395").Value = 45180
$ws.Range("E395").Value = 16
$ws.Range("F395").Value = 100112017
$ws.Range("G395").Value = "Apio"
$ws.Range("H395").Value = "Americana (o)"
$ws.Range("I395").Value = "Primera"
$ws.Range("J395").Value = 120
$ws.Range("K395").Value = 6000
$ws.Range("L395").Value = 6500
$ws.Range("M395").Value = 6250
$ws.Range("N395").Value = "$/docena de matas"
$ws.Range("O395").Value = "Provincia del Elquí"
$ws.Range("P395").Value = 1042
$ws.Range("Q395").Value = 6
$ws.Range("R395").Value = "Hortaliza"
